$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    2  = -1
    5  = 2
    8  = 1
    11 = 4
    13 = 1
    14 = 3
    15 = 0
    17 = 0
    20 = 1
    24 = 1
    27 = -1
    29 = -2
    35 = 1
    36 = 2
    37 = 3
    43 = 5
    57 = 0
    61 = -4
    62 = -3
    65 = -3
    66 = -4
    67 = -1
    68 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
